$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right, 5-twip space) to the first
# paragraph, matching the style already used elsewhere in the document.
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromRight = 5

# Increase the left indent from 120 to 225 twips (11.25 pt).
$p1.LeftIndent = 11.25

# The paragraph currently holds two runs: the placeholder text, then a
# trailing run containing a single space. The paragraph's Range ends one
# position past the paragraph mark, so the space sits immediately before
# that, and the text run sits immediately before the space.
$pStart = $p1.Range.Start
$pEnd = $p1.Range.End

# Drop the trailing run that contains only a space.
$d.Range($pEnd - 2, $pEnd - 1).Delete()

# Update the placeholder text in the (now sole) run of the paragraph.
$d.Range($pStart, $pEnd - 2).Text = "**ID__AFFARS_5342_7100__ID**"
